$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the "researchers per million" table with the next data
# year (2020 / 534), matching the styling already used for the other
# year-header and value cells in columns D:M. -----------------------

# N4: year header, same formatting as the other plain year headers
# (D4:K4 use style 12; L4/M4 use style 13 -- the header row keeps the
# plain right-aligned look, so copy from D4).
$ws.Range("D4").Copy()
$ws.Range("N4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("N4").Value = 2020

# N5: data value, matching the bordered style used by the two most
# recent years (L5:M5 -> style 17).
$ws.Range("L5").Copy()
$ws.Range("N5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("N5").Value = 534

$excel.CutCopyMode = 0

# --- Update the saved view state: scroll the window so column E is
# the left-most visible column, and leave the active selection on
# S10. -----------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("S10").Select() | Out-Null
